$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (serial date 45406 -> 45436, i.e. 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = (Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0).Date

# Update the price list in column D (rows 30-37) with the new prices
$ws.Range("D30").Value = 120.069
$ws.Range("D31").Value = 128.588
$ws.Range("D32").Value = 140.928
$ws.Range("D33").Value = 208.099
$ws.Range("D34").Value = 269.256
$ws.Range("D35").Value = 379.682
$ws.Range("D36").Value = 550.922
$ws.Range("D37").Value = 844.688
